$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# Update formulas in row 4: change constant 50 -> 235
$ws.Range("A4").Formula = "=1/(235*0.0002*A2^-27.15)"
$ws.Range("B4:J4").FormulaR1C1 = "=1/(235*0.0002*R[-2]C^-27.15)"

# Update the selection on the sheet to the full data range A1:J4
$ws.Activate()
$ws.Range("A1:J4").Select()
